$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for rows 2-8, columns A (TruckID), B (AssignedDockPosition),
# C (start_loading_time), D (end_loading_time)
$data = @{
    2 = @(3, 1, 5, 7)
    3 = @(4, 1, 12, 12)
    4 = @(6, 1, 17, 17)
    5 = @(7, 1, 22, 22)
    6 = @(5, 2, 5, 5)
    7 = @(1, 3, 5, 5)
    8 = @(2, 3, 10, 11)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 1).Value = $values[0]
    $ws.Cells.Item($row, 2).Value = $values[1]
    $ws.Cells.Item($row, 3).Value = $values[2]
    $ws.Cells.Item($row, 4).Value = $values[3]
}
